# "Updated in class 5"
#
# Slide 14 has a Cylinder class diagram. The fields box (group "组合 4",
# shape "矩形 15") lists three fields:
#   pi: double
#   Radius: double
#   Height: int
#
# Fix the capitalization of the last two field names so they follow the
# same lower-case Java field-naming convention as "pi":
#   Radius: double -> radius: double
#   Height: int    -> height: int

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(14)
$group = $slide.Shapes.Item("组合 4")
$fieldsShape = $group.GroupItems.Item("矩形 15")
$tr = $fieldsShape.TextFrame.TextRange

# Full text is "pi: double" + CR + "Radius: double" + CR + "Height: int"
# so (1-based) character offsets are:
#   1-10  : "pi: double"
#   11    : paragraph break
#   12-25 : "Radius: double"   (12 = "R", 13-17 = "adius", 18-25 = ": double")
#   26    : paragraph break
#   27-37 : "Height: int"

# Retype the leading "R" as "r" ...
$tr.Characters(12, 1).Text = "r"

# ... then retype the trailing ": double" in place, which splits the
# remainder of the word ("adius") into its own run, matching how the
# correction was actually typed in PowerPoint.
$tr.Characters(18, 8).Text = ": double"

# "Height: int" -> "height: int"
$tr.Characters(27, 11).Text = "height: int"
